$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 (shifts existing rows 4-8 down to 5-9)
$ws.Rows("4:4").Insert()

# Populate the newly inserted row 4 with the new record's data
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C4").Value = "Coquimbo"
$ws.Range("D4").Value = 44497
$ws.Range("D4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100101
$ws.Range("H4").Value = "Berries"
$ws.Range("I4").Value = 100101001
$ws.Range("J4").Value = "Arándano (blue)"
$ws.Range("K4").Value = "Sin especificar"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 500
$ws.Range("N4").Value = 9000
$ws.Range("O4").Value = 10000
$ws.Range("P4").Value = 9500
$ws.Range("Q4").Value = "`$/bandeja 2 kilos"
$ws.Range("R4").Value = "Provincia de Limarí"
$ws.Range("S4").Value = 4750
$ws.Range("T4").Value = 2
